# Fix a bug on the generated "inputs" rows of the Squelette_sujet_GUI sheet:
#  - regenerate the 4 stimulus rows (values + timings) that were produced with
#    a faulty randomizer (duplicate/garbage "ERREUR" entries in the Erreur column)
#  - widen column A (Stimulus) and column F (TR) so the new content fits

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668   # -> stored width 21
$ws.Columns.Item(6).ColumnWidth = 8.877604166666666    # -> stored width ~9.71

# --- regenerated data rows -------------------------------------------------
# Columns: A Stimulus | B Déterminant_Mot | C Nom_Mot | D Déterminant_image
#          E Nom_Image | F TR | G Lettre | H Congruence | I Erreur

$rows = @(
    @{ Row = 2;  A = "Squelette_sujet_gui_1"; B = "La";  C = "pomme";    D = "Des"; E = "pommes";    F = 0.37934820004738867; G = "e"; H = $false; I = 0 },
    @{ Row = 3;  A = "Squelette_sujet_gui_2"; B = "Des"; C = "carottes"; D = "Une"; E = "carotte";    F = 1.1503209000220522;  G = "e"; H = $false; I = 0 },
    @{ Row = 4;  A = "Squelette_sujet_gui_3"; B = "Une"; C = "tomate";   D = "Les"; E = "tomates";    F = 0.38516549998894334; G = "e"; H = $false; I = 0 },
    @{ Row = 5;  A = "Squelette_sujet_gui_4"; B = "Le";  C = "broccoli"; D = "Un";  E = "broccoli";   F = 0.24058650003280491; G = "q"; H = $true;  I = 0 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H

    # Column I ("Erreur") used to contain the literal text "ERREUR" on the
    # broken rows -- that was the bug. Force a real numeric 0 even though the
    # column is still formatted as text, instead of the stray error flag.
    $cell = $ws.Cells.Item($n, 9)
    $cell.NumberFormat = "General"
    $cell.Value = $r.I
    $cell.NumberFormat = "@"
}
